$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.256161925124957
$ws.Range("C2").Value = 0.2346248611648036
$ws.Range("D2").Value = 0.07810553557496291
$ws.Range("E2").Value = 0.104616351433313
$ws.Range("G2").Value = 0.002472023795564883
$ws.Range("L2").Value = 0.1906439604913075
$ws.Range("M2").Value = 0.2632552252174563
$ws.Range("O2").Value = 4.235721211673194
$ws.Range("B3").Value = 1.156852820408346
$ws.Range("C3").Value = 0.2197061655329833
$ws.Range("D3").Value = 0.07090213355218111
$ws.Range("E3").Value = 0.1055377259485465
$ws.Range("G3").Value = 0.002475713788156935
$ws.Range("L3").Value = 0.1881033970932435
$ws.Range("M3").Value = 0.248139918629434
$ws.Range("O3").Value = 4.240221910131282
$ws.Range("B4").Value = 1.096187617494081
$ws.Range("C4").Value = 0.2104786647053061
$ws.Range("D4").Value = 0.06651650428111111
$ws.Range("E4").Value = 0.1061351793882015
$ws.Range("G4").Value = 0.002478100188064472
$ws.Range("L4").Value = 0.186637737513216
$ws.Range("M4").Value = 0.2389541189392261
$ws.Range("O4").Value = 4.245968564097893
$ws.Range("B5").Value = 1.071545227969978
$ws.Range("C5").Value = 0.2067016074725814
$ws.Range("D5").Value = 0.06473863993083739
$ws.Range("E5").Value = 0.1063866352511141
$ws.Range("G5").Value = 0.002479103119353699
$ws.Range("L5").Value = 0.1860641951582522
$ws.Range("M5").Value = 0.2352348752442737
$ws.Range("O5").Value = 4.249058840838387
$ws.Range("B6").Value = 1.067458188885723
$ws.Range("C6").Value = 0.2060734216453
$ws.Range("D6").Value = 0.06444398774344506
$ws.Range("E6").Value = 0.106428872130132
$ws.Range("G6").Value = 0.002479271497478326
$ws.Range("L6").Value = 0.185970392597568
$ws.Range("M6").Value = 0.2346187540856306
$ws.Range("O6").Value = 4.249617135365639
$ws.Range("B7").Value = 1.095854959739199
$ws.Range("C7").Value = 0.2104277936965104
$ws.Range("D7").Value = 0.06649248974386524
$ws.Range("E7").Value = 0.1061385382490986
$ws.Range("G7").Value = 0.002478113590480446
$ws.Range("L7").Value = 0.1866299064226666
$ws.Range("M7").Value = 0.2389038623750466
$ws.Range("O7").Value = 4.246007212382636
$ws.Range("B8").Value = 1.221856145922629
$ws.Range("C8").Value = 0.2294949682484173
$ws.Range("D8").Value = 0.07561402404675732
$ws.Range("E8").Value = 0.1049274633929027
$ws.Range("G8").Value = 0.002473271105165196
$ws.Range("L8").Value = 0.1897484243759493
$ws.Range("M8").Value = 0.2580238163529103
$ws.Range("O8").Value = 4.236652863790624
$ws.Range("B9").Value = 1.471382804183861
$ws.Range("C9").Value = 0.2663458279110671
$ws.Range("D9").Value = 0.09380091705030225
$ws.Range("E9").Value = 0.1028037576126187
$ws.Range("G9").Value = 0.002464728573101022
$ws.Range("L9").Value = 0.1966112207210813
$ws.Range("M9").Value = 0.2962682912531349
$ws.Range("O9").Value = 4.242060850651342
$ws.Range("B10").Value = 1.656178493525317
$ws.Range("C10").Value = 0.2930865997379897
$ws.Range("D10").Value = 0.1073519963765079
$ws.Range("E10").Value = 0.1013958593575999
$ws.Range("G10").Value = 0.00245902757286694
$ws.Range("L10").Value = 0.2021091352093833
$ws.Range("M10").Value = 0.3248218966814918
$ws.Range("O10").Value = 4.260634548465845
$ws.Range("B11").Value = 1.740563399697692
$ws.Range("C11").Value = 0.3051785082445804
$ws.Range("D11").Value = 0.1135592973275266
$ws.Range("E11").Value = 0.1007882982870888
$ws.Range("G11").Value = 0.002456557635080697
$ws.Range("L11").Value = 0.2047093603541299
$ws.Range("M11").Value = 0.3379103786546906
$ws.Range("O11").Value = 4.272281078242713
$ws.Range("B12").Value = 1.772563234300208
$ws.Range("C12").Value = 0.3097468440111015
$ws.Range("D12").Value = 0.1159160915609192
$ws.Range("E12").Value = 0.1005629501795818
$ws.Range("G12").Value = 0.002455639990721154
$ws.Range("L12").Value = 0.2057082544896787
$ws.Range("M12").Value = 0.3428808442054319
$ws.Range("O12").Value = 4.277153087318936
$ws.Range("B13").Value = 1.76566950674794
$ws.Range("C13").Value = 0.3087634455687294
$ws.Range("D13").Value = 0.1154082359359165
$ws.Range("E13").Value = 0.1006112730488561
$ws.Range("G13").Value = 0.002455836837462866
$ws.Range("L13").Value = 0.2054924914979779
$ws.Range("M13").Value = 0.3418097387408849
$ws.Range("O13").Value = 4.276083246605424
$ws.Range("B14").Value = 1.743195146949063
$ws.Range("C14").Value = 0.3055545615149811
$ws.Range("D14").Value = 0.1137530671067992
$ws.Range("E14").Value = 0.1007696641748427
$ws.Range("G14").Value = 0.002456481786339292
$ws.Range("L14").Value = 0.2047912545853734
$ws.Range("M14").Value = 0.3383190192123848
$ws.Range("O14").Value = 4.27267263655807
$ws.Range("B15").Value = 1.729434797671217
$ws.Range("C15").Value = 0.3035876410520757
$ws.Range("D15").Value = 0.1127400408327759
$ws.Range("E15").Value = 0.1008672980097693
$ws.Range("G15").Value = 0.002456879134488276
$ws.Range("L15").Value = 0.2043635813341496
$ws.Range("M15").Value = 0.3361826902927447
$ws.Range("O15").Value = 4.270643727300069
$ws.Range("B16").Value = 1.650670105429015
$ws.Range("C16").Value = 0.2922948897420383
$ws.Range("D16").Value = 0.1069472022103355
$ws.Range("E16").Value = 0.1014362261635195
$ws.Range("G16").Value = 0.002459191466179742
$ws.Range("L16").Value = 0.2019411988265176
$ws.Range("M16").Value = 0.3239685199628966
$ws.Range("O16").Value = 4.259937924811425
$ws.Range("B17").Value = 1.602431886155216
$ws.Range("C17").Value = 0.2853484371215416
$ws.Range("D17").Value = 0.1034045020651888
$ws.Range("E17").Value = 0.1017936651123303
$ws.Range("G17").Value = 0.0024606415690749
$ws.Range("L17").Value = 0.2004805384645039
$ws.Range("M17").Value = 0.316500856853331
$ws.Range("O17").Value = 4.25419055009317
$ws.Range("B18").Value = 1.574716754767451
$ws.Range("C18").Value = 0.2813461973832716
$ws.Range("D18").Value = 0.1013708638482171
$ws.Range("E18").Value = 0.1020023516595828
$ws.Range("G18").Value = 0.002461487257214697
$ws.Range("L18").Value = 0.1996497431811122
$ws.Range("M18").Value = 0.3122150095606671
$ws.Range("O18").Value = 4.251185636258953
$ws.Range("B19").Value = 1.565338107191792
$ws.Range("C19").Value = 0.2799899409155273
$ws.Range("D19").Value = 0.1006829990001279
$ws.Range("E19").Value = 0.1020735415395126
$ws.Range("G19").Value = 0.002461775592536295
$ws.Range("L19").Value = 0.1993700545586705
$ws.Range("M19").Value = 0.3107655063003847
$ws.Range("O19").Value = 4.250219831278088
$ws.Range("B20").Value = 1.607563806771452
$ws.Range("C20").Value = 0.2860886063356531
$ws.Range("D20").Value = 0.1037812110834579
$ws.Range("E20").Value = 0.1017552946495912
$ws.Range("G20").Value = 0.002460486000435309
$ws.Range("L20").Value = 0.2006350619531929
$ws.Range("M20").Value = 0.3172948351417659
$ws.Range("O20").Value = 4.25477121957627
$ws.Range("B21").Value = 1.749795194819967
$ws.Range("C21").Value = 0.3064973777497073
$ws.Range("D21").Value = 0.1142390608936807
$ws.Range("E21").Value = 0.1007230127748369
$ws.Range("G21").Value = 0.002456291870506337
$ws.Range("L21").Value = 0.2049968384000067
$ws.Range("M21").Value = 0.3393439454461102
$ws.Range("O21").Value = 4.273661868804425
$ws.Range("B22").Value = 1.8430140796508
$ws.Range("C22").Value = 0.3197737854236777
$ws.Range("D22").Value = 0.1211101782462407
$ws.Range("E22").Value = 0.1000758768224305
$ws.Range("G22").Value = 0.002453653702398921
$ws.Range("L22").Value = 0.2079305355896679
$ws.Range("M22").Value = 0.3538366788954121
$ws.Range("O22").Value = 4.288700174379187
$ws.Range("B23").Value = 1.793237709069444
$ws.Range("C23").Value = 0.312693635451609
$ws.Range("D23").Value = 0.1174395906385683
$ws.Range("E23").Value = 0.1004187502489493
$ws.Range("G23").Value = 0.002455052352544683
$ws.Range("L23").Value = 0.2063571755994502
$ws.Range("M23").Value = 0.3460941428186715
$ws.Range("O23").Value = 4.280426965737092
$ws.Range("B24").Value = 1.605243610371758
$ws.Range("C24").Value = 0.2857540027114283
$ws.Range("D24").Value = 0.1036108912651201
$ws.Range("E24").Value = 0.1017726319980391
$ws.Range("G24").Value = 0.002460556295506987
$ws.Range("L24").Value = 0.2005651739806638
$ws.Range("M24").Value = 0.3169358544614695
$ws.Range("O24").Value = 4.254507766628137
$ws.Range("B25").Value = 1.403620110851762
$ws.Range("C25").Value = 0.2564349597737987
$ws.Range("D25").Value = 0.08884813298874406
$ws.Range("E25").Value = 0.1033514610825135
$ws.Range("G25").Value = 0.002466938103495412
$ws.Range("L25").Value = 0.1946746217787236
$ws.Range("M25").Value = 0.3418097387408849
$ws.Range("O25").Value = 4.23804258059198
